$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("MOGOLLON MORON KARIN LISSET", 47),
    @("CHIROQUE YARLEQUE BETTY ELIZABETH", 39),
    @("TUANAMA PIZANGO ELIZABETH", 38),
    @("SILVA ALVARADO EVELYN DE JESUS", 38),
    @("AGURTO TINEO CESIA JIMENA", 37),
    @("MORETO ESPINOZA CRISTIAN ESTEBAN", 36),
    @("TENE TRABUCCO GIAN PIERRE", 36),
    @("NEYRA PEREIRA MONICA AGNES ALEXANDRA", 35),
    @("SANCARRANCO SANCHEZ DE CRUZ GISSELA SHANI", 33),
    @("HERRERA JUAN MANUEL", 31)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
